$d = $word.ActiveDocument

$pairs = @(
  @{old='23-5='; new='35-24='},
  @{old='41-14='; new='44+0='},
  @{old='53-15='; new='12+73='},
  @{old='0+44='; new='49-10='},
  @{old='86+5='; new='46-5='},
  @{old='86-57='; new='70-7='},
  @{old='36-34='; new='54+13='},
  @{old='82-75='; new='32-13='},
  @{old='39+35='; new='62+26='},
  @{old='56+41='; new='21+22='},
  @{old='72-21='; new='78-71='},
  @{old='64-30='; new='82-76='},
  @{old='7+43='; new='20-0='},
  @{old='14+8='; new='73-67='},
  @{old='51-16='; new='76-30='},
  @{old='95-81='; new='54+10='},
  @{old='9+24='; new='14+76='},
  @{old='6-1='; new='8+30='},
  @{old='40-14='; new='71-63='},
  @{old='35+30='; new='17+35='},
  @{old='28+38='; new='64-57='},
  @{old='39+33='; new='58-50='},
  @{old='46+23='; new='29-27='},
  @{old='35-23='; new='14+13='},
  @{old='63+11='; new='38-25='},
  @{old='19+27='; new='8+18='},
  @{old='61-34='; new='95-60='},
  @{old='35-1='; new='81-22='},
  @{old='98-3='; new='3+83='},
  @{old='92-10='; new='51-14='},
  @{old='77-21='; new='2+73='},
  @{old='5+14='; new='82-33='},
  @{old='15+70='; new='79-72='},
  @{old='96-68='; new='53+46='},
  @{old='71-60='; new='35+38='},
  @{old='27-11='; new='37-15='},
  @{old='23+26='; new='92-29='},
  @{old='94-68='; new='79-58='},
  @{old='26+22='; new='61-57='},
  @{old='42+13='; new='76-65='},
  @{old='42-4='; new='54-38='},
  @{old='80-50='; new='55-27='},
  @{old='90-48='; new='5+54='},
  @{old='38-21='; new='89-72='},
  @{old='0+38='; new='14+43='},
  @{old='63-2='; new='22-6='},
  @{old='62-32='; new='62-39='},
  @{old='64-21='; new='63+27='},
  @{old='35+2='; new='8+48='},
  @{old='82-69='; new='10+63='},
  @{old='60-11='; new='40-20='},
  @{old='82-18='; new='57+5='},
  @{old='37+9='; new='53+30='},
  @{old='9+82='; new='84-5='},
  @{old='22-1='; new='99-15='},
  @{old='75-17='; new='88-71='},
  @{old='13+30='; new='52-36='},
  @{old='20+70='; new='80-25='},
  @{old='22+39='; new='18-0='},
  @{old='64+10='; new='79-51='},
  @{old='31-31='; new='90-49='},
  @{old='36-15='; new='84-20='},
  @{old='42+7='; new='60-45='},
  @{old='62-33='; new='54-3='},
  @{old='63-39='; new='83-58='},
  @{old='62-36='; new='38+4='},
  @{old='36+34='; new='3+62='},
  @{old='49-22='; new='76-9='},
  @{old='50-44='; new='56-11='},
  @{old='57+13='; new='4+6='},
  @{old='52-51='; new='10+16='},
  @{old='69-67='; new='99-46='},
  @{old='2+10='; new='21+75='},
  @{old='5+68='; new='4-1='},
  @{old='62-3='; new='49-35='},
  @{old='98-89='; new='32+42='},
  @{old='62-22='; new='22+30='},
  @{old='94-10='; new='20+71='},
  @{old='38+13='; new='43-40='},
  @{old='50-40='; new='21+21='},
  @{old='98-78='; new='89-43='},
  @{old='51+1='; new='41+48='},
  @{old='20+56='; new='21+8='},
  @{old='36+52='; new='95-26='},
  @{old='72-57='; new='12+39='},
  @{old='71-37='; new='88-47='},
  @{old='34+35='; new='5+52='},
  @{old='47-45='; new='50-14='},
  @{old='63+34='; new='89-8='},
  @{old='60+0='; new='54+35='},
  @{old='95-61='; new='40-26='},
  @{old='74+6='; new='33-14='},
  @{old='97-94='; new='25+18='},
  @{old='8+50='; new='20+35='},
  @{old='50-47='; new='76-59='},
  @{old='67-19='; new='75-66='},
  @{old='41+47='; new='87-62='},
  @{old='96-15='; new='63+36='},
  @{old='10+23='; new='56-42='},
  @{old='0+94='; new='70+26='},
)

foreach ($p in $pairs) {
  $d.Content.Find.Execute($p.old, $true, $true, $false, $false, $false, $true, 1, $false, $p.new, 2) | Out-Null
}
